# Finished writeup and minor edits.
#
# The document currently ends with a single empty bulleted (ListParagraph)
# paragraph that only contains the hidden "_GoBack" bookmark. We need to
# turn that into three bulleted paragraphs of write-up text, each with
# double line spacing, while keeping the "_GoBack" bookmark anchored in the
# last paragraph (immediately before the final trailing space run), exactly
# as it was before (just shifted along with the new content around it).

$d = $word.ActiveDocument

# The final paragraph in the document is the empty bulleted paragraph that
# holds the "_GoBack" bookmark.
$lastIndex = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($lastIndex)
$r = $p.Range
$r.Collapse(1)                 # wdCollapseStart
$r.InsertParagraphBefore() | Out-Null
$r.InsertParagraphBefore() | Out-Null

# We now have three consecutive bulleted paragraphs; the bookmark is still
# anchored in the last (third) one.
$p1 = $d.Paragraphs.Item($lastIndex)
$p2 = $d.Paragraphs.Item($lastIndex + 1)
$p3 = $d.Paragraphs.Item($lastIndex + 2)

# Insert the lone trailing space run *after* the bookmark first, while the
# bookmark is still the unique thing sitting at the very end of the
# document -- this is what makes the new run land after it rather than
# merging into the text that will be placed before it in the next step.
$trailing = $d.Range($d.Content.End, $d.Content.End)
$trailing.InsertAfter(" ")

# Bullet 1
$p1.Range.InsertBefore('There is a correlation between the most purchased items and the most profitable items. This should be no surprise, as the most profitable items are calculated using the amount purchased. For dataset 2, the same items appear on both the top five lists of most purchased and most profitable items. However, for the first dataset, this is not the case. There was only one item (the Retribution Axe) that was on both top five lists. This could be because all the other most purchased items cost less than $2.50 and were likely purchased more often due to the lower price. ')

# Bullet 2
$p2.Range.InsertBefore('Overall, there are vastly more male players than female players, around 81% male versus 17.6% female with the remainder made up of other/non-disclosed players. This explains the greater overall purchase counts by males. However, in both datasets there are also more repeat purchases by male players. Perhaps the genre of the game appeals more to males or perhaps the marketing heavily targeted them.')

# Bullet 3 (text goes before the bookmark; the trailing space run we already
# created stays after the bookmark).
$p3.Range.InsertBefore('For both datasets, there appears to be an approximate normal distribution for player age on purchase count. There are big spikes at the 20-24 age range with counts in all the surrounding age ranges. This age range did not stand apart from the others in average purchasing prices. Their increased profitability came solely due to their large number of purchases. It could be that people in this age range tend to be more attracted to this genre of games or to gaming in general. It could also be that marketing that heavily targeted this age group was effective. ')

# Apply double line spacing (480 twips / auto rule) to all three bullets.
$p1.LineSpacingRule = 2   # wdLineSpaceDouble
$p2.LineSpacingRule = 2   # wdLineSpaceDouble
$p3.LineSpacingRule = 2   # wdLineSpaceDouble
